$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.288150666666667
$ws.Range("H2").Value = 6.864452
$ws.Range("I2").Value = 0.3964219041944151
$ws.Range("J2").Value = 0.3964219041944151
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.793878
$ws.Range("N2").Value = 32.381634
$ws.Range("O2").Value = 0.166691744666625
$ws.Range("P2").Value = 0.166691744666625
$ws.Range("Q2").Value = 24.69801914161867
$ws.Range("R2").Value = 222.282172274568
$ws.Range("S2").Value = 0.0660802588342327
$ws.Range("T2").Value = 0.0660802588342327

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.288150666666667
$ws.Range("H3").Value = 6.864452
$ws.Range("I3").Value = 0.3964219041944151
$ws.Range("J3").Value = 0.3964219041944151
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.04277733333333
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4485126868821051
$ws.Range("P3").Value = 0.4485126868821051
$ws.Range("Q3").Value = 66.45425031711822
$ws.Range("R3").Value = 598.088252854064
$ws.Range("S3").Value = 0.1778002533891576
$ws.Range("T3").Value = 0.1778002533891576

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.288150666666667
$ws.Range("H4").Value = 6.864452
$ws.Range("I4").Value = 0.3964219041944151
$ws.Range("J4").Value = 0.3964219041944151
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.345855666666666
$ws.Range("N4").Value = 16.037567
$ws.Range("O4").Value = 0.08255698348755007
$ws.Range("P4").Value = 0.08255698348755008
$ws.Range("Q4").Value = 12.23212320758711
$ws.Range("R4").Value = 110.089108868284
$ws.Range("S4").Value = 0.03272739659868149
$ws.Range("T4").Value = 0.03272739659868149

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.288150666666667
$ws.Range("H5").Value = 6.864452
$ws.Range("I5").Value = 0.3964219041944151
$ws.Range("J5").Value = 0.3964219041944151
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.57101366666667
$ws.Range("N5").Value = 58.713041
$ws.Range("O5").Value = 0.3022385849637199
$ws.Range("P5").Value = 0.3022385849637199
$ws.Range("Q5").Value = 44.78142796872577
$ws.Range("R5").Value = 403.032851718532
$ws.Range("S5").Value = 0.1198139953723434
$ws.Range("T5").Value = 0.1198139953723434

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.588894
$ws.Range("H6").Value = 7.766681999999999
$ws.Range("I6").Value = 0.4485256605643812
$ws.Range("J6").Value = 0.4485256605643813
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.793878
$ws.Range("N6").Value = 32.381634
$ws.Range("O6").Value = 0.166691744666625
$ws.Range("P6").Value = 0.166691744666625
$ws.Range("Q6").Value = 27.944205990932
$ws.Range("R6").Value = 251.497853918388
$ws.Range("S6").Value = 0.07476552488722712
$ws.Range("T6").Value = 0.07476552488722714

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.588894
$ws.Range("H7").Value = 7.766681999999999
$ws.Range("I7").Value = 0.4485256605643812
$ws.Range("J7").Value = 0.4485256605643813
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 29.04277733333333
$ws.Range("N7").Value = 87.128332
$ws.Range("O7").Value = 0.4485126868821051
$ws.Range("P7").Value = 0.4485126868821051
$ws.Range("Q7").Value = 75.18867198160267
$ws.Range("R7").Value = 676.698047834424
$ws.Range("S7").Value = 0.2011694491553017
$ws.Range("T7").Value = 0.2011694491553017

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.588894
$ws.Range("H8").Value = 7.766681999999999
$ws.Range("I8").Value = 0.4485256605643812
$ws.Range("J8").Value = 0.4485256605643813
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.345855666666666
$ws.Range("N8").Value = 16.037567
$ws.Range("O8").Value = 0.08255698348755007
$ws.Range("P8").Value = 0.08255698348755008
$ws.Range("Q8").Value = 13.83985366029933
$ws.Range("R8").Value = 124.558682942694
$ws.Range("S8").Value = 0.03702892555295611
$ws.Range("T8").Value = 0.03702892555295612

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.588894
$ws.Range("H9").Value = 7.766681999999999
$ws.Range("I9").Value = 0.4485256605643812
$ws.Range("J9").Value = 0.4485256605643813
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.57101366666667
$ws.Range("N9").Value = 58.713041
$ws.Range("O9").Value = 0.3022385849637199
$ws.Range("P9").Value = 0.3022385849637199
$ws.Range("Q9").Value = 50.66727985555133
$ws.Range("R9").Value = 456.005518699962
$ws.Range("S9").Value = 0.1355617609688963
$ws.Range("T9").Value = 0.1355617609688964

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3123523333333333
$ws.Range("H10").Value = 0.9370569999999999
$ws.Range("I10").Value = 0.05411501461132016
$ws.Range("J10").Value = 0.05411501461132018
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.793878
$ws.Range("N10").Value = 32.381634
$ws.Range("O10").Value = 0.166691744666625
$ws.Range("P10").Value = 0.166691744666625
$ws.Range("Q10").Value = 3.371492979015333
$ws.Range("R10").Value = 30.34343681113799
$ws.Range("S10").Value = 0.009020526198220859
$ws.Range("T10").Value = 0.009020526198220861

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3123523333333333
$ws.Range("H11").Value = 0.9370569999999999
$ws.Range("I11").Value = 0.05411501461132016
$ws.Range("J11").Value = 0.05411501461132018
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 29.04277733333333
$ws.Range("N11").Value = 87.128332
$ws.Range("O11").Value = 0.4485126868821051
$ws.Range("P11").Value = 0.4485126868821051
$ws.Range("Q11").Value = 9.07157926654711
$ws.Range("R11").Value = 81.64421339892399
$ws.Range("S11").Value = 0.02427127060398758
$ws.Range("T11").Value = 0.02427127060398759

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.3123523333333333
$ws.Range("H12").Value = 0.9370569999999999
$ws.Range("I12").Value = 0.05411501461132016
$ws.Range("J12").Value = 0.05411501461132018
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 5.345855666666666
$ws.Range("N12").Value = 16.037567
$ws.Range("O12").Value = 0.08255698348755007
$ws.Range("P12").Value = 0.08255698348755008
$ws.Range("Q12").Value = 1.669790491146555
$ws.Range("R12").Value = 15.028114420319
$ws.Range("S12").Value = 0.00446757236769529
$ws.Range("T12").Value = 0.004467572367695291

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.3123523333333333
$ws.Range("H13").Value = 0.9370569999999999
$ws.Range("I13").Value = 0.05411501461132016
$ws.Range("J13").Value = 0.05411501461132018
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.57101366666667
$ws.Range("N13").Value = 58.713041
$ws.Range("O13").Value = 0.3022385849637199
$ws.Range("P13").Value = 0.3022385849637199
$ws.Range("Q13").Value = 6.113051784481888
$ws.Range("R13").Value = 55.01746606033699
$ws.Range("S13").Value = 0.01635564544141643
$ws.Range("T13").Value = 0.01635564544141644

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5826116666666666
$ws.Range("H14").Value = 1.747835
$ws.Range("I14").Value = 0.1009374206298835
$ws.Range("J14").Value = 0.1009374206298836
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 10.793878
$ws.Range("N14").Value = 32.381634
$ws.Range("O14").Value = 0.166691744666625
$ws.Range("P14").Value = 0.166691744666625
$ws.Range("Q14").Value = 6.288639251376666
$ws.Range("R14").Value = 56.59775326238999
$ws.Range("S14").Value = 0.01682543474694427
$ws.Range("T14").Value = 0.01682543474694427

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5826116666666666
$ws.Range("H15").Value = 1.747835
$ws.Range("I15").Value = 0.1009374206298835
$ws.Range("J15").Value = 0.1009374206298836
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 29.04277733333333
$ws.Range("N15").Value = 87.128332
$ws.Range("O15").Value = 0.4485126868821051
$ws.Range("P15").Value = 0.4485126868821051
$ws.Range("Q15").Value = 16.92066090680222
$ws.Range("R15").Value = 152.28594816122
$ws.Range("S15").Value = 0.0452717137336583
$ws.Range("T15").Value = 0.0452717137336583

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5826116666666666
$ws.Range("H16").Value = 1.747835
$ws.Range("I16").Value = 0.1009374206298835
$ws.Range("J16").Value = 0.1009374206298836
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.345855666666666
$ws.Range("N16").Value = 16.037567
$ws.Range("O16").Value = 0.08255698348755007
$ws.Range("P16").Value = 0.08255698348755008
$ws.Range("Q16").Value = 3.114557879716111
$ws.Range("R16").Value = 28.03102091744499
$ws.Range("S16").Value = 0.008333088968217191
$ws.Range("T16").Value = 0.008333088968217195

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5826116666666666
$ws.Range("H17").Value = 1.747835
$ws.Range("I17").Value = 0.1009374206298835
$ws.Range("J17").Value = 0.1009374206298836
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.57101366666667
$ws.Range("N17").Value = 58.713041
$ws.Range("O17").Value = 0.3022385849637199
$ws.Range("P17").Value = 0.3022385849637199
$ws.Range("Q17").Value = 11.40230089069278
$ws.Range("R17").Value = 102.620708016235
$ws.Range("S17").Value = 0.03050718318106379
$ws.Range("T17").Value = 0.0305071831810638
